# Update the dSF (column F) values to reflect the repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -1
    4  = -2
    6  = -8
    7  = -4
    8  = -1
    9  = -5
    10 = 4
    11 = -1
    12 = -6
    13 = -4
    15 = -4
    16 = -4
    17 = 1
    18 = -3
    19 = 2
    20 = -4
    21 = -1
    22 = 2
    23 = -3
    25 = 2
    26 = -6
    28 = 0
    29 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
